$d = $word.ActiveDocument
$wNs = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# ---------------------------------------------------------------------------
# Edit 1: first paragraph - bold "Çalışanlar" (capitalised, split Ç/alışanlar)
# ---------------------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$xml1 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00930A40" w:rsidRDefault="007E3FF3"><w:r><w:t xml:space="preserve">Aşağıdaki özelliklere sahip </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>Ç</w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t>alışanlar</w:t></w:r><w:r><w:t xml:space="preserve"> tablosu oluşturun ve 5 adet kayıt ekleyin(MS Access kullanarak)</w:t></w:r></w:p>
'@
$p1.InsertXML($xml1)

# ---------------------------------------------------------------------------
# Edit 2: "SELECT * FROM Çalışanlar WHERE DoğumYeri = "İstanbul""
# ---------------------------------------------------------------------------
$p2 = $d.Paragraphs(12).Range
$xml2 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve">SELECT * FROM Çalışanlar WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009827A8"><w:t>DoğumYeri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009827A8"><w:t>="İstanbul"</w:t></w:r></w:p>
'@
$p2.InsertXML($xml2)

# ---------------------------------------------------------------------------
# Edit 3: "SELECT * FROM Çalışanlar WHERE Birimi = "Muhasebe""
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(14).Range
$xml3 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t>SELECT * FROM Çalışanlar WHERE Birimi="Muhasebe"</w:t></w:r></w:p>
'@
$p3.InsertXML($xml3)

# ---------------------------------------------------------------------------
# Edit 4: Sorgu 3 SELECT (TCKimlik,Adı,Soyadı FROM Çalışanlar WHERE Maaşı>2500)
# ---------------------------------------------------------------------------
$p4 = $d.Paragraphs(16).Range
$xml4 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="009827A8"><w:t>TCKimlik,Adı</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="009827A8"><w:t>,Soyadı</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve"> FROM Çalışanlar WHERE Maaşı&gt;2500</w:t></w:r></w:p>
'@
$p4.InsertXML($xml4)

# ---------------------------------------------------------------------------
# Edit 5: Sorgu 4 SELECT (Adı,Soyadı,DoğumYeri FROM Çalışanlar ORDER BY Maaşı DESC)
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs(18).Range
$xml5 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="009827A8"><w:t>Adı,Soyadı</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="009827A8"><w:t>,DoğumYeri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve"> FROM Çalışanlar ORDER BY Maaşı DESC</w:t></w:r></w:p>
'@
$p5.InsertXML($xml5)

# ---------------------------------------------------------------------------
# Edit 6: Sorgu 5 SELECT (Adı,Soyadı,Maaşı FROM Çalışanlar WHERE DoğumYeri<>"Denizli"  ORDER BY Maaşı)
# ---------------------------------------------------------------------------
$p6 = $d.Paragraphs(20).Range
$xml6 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve">SELECT </w:t></w:r><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidRPr="009827A8"><w:t>Adı,Soyadı</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidRPr="009827A8"><w:t>,Maaşı</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009827A8"><w:t xml:space="preserve"> FROM Çalışanlar WHERE </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="009827A8"><w:t>DoğumYeri</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="009827A8"><w:t>&lt;&gt;"Denizli"  ORDER BY Maaşı</w:t></w:r></w:p>
'@
$p6.InsertXML($xml6)

# ---------------------------------------------------------------------------
# Edit 7: Sorgu 6 SELECT (SELECT Maaşı FROM Çalışanlar WHERE Birimi IS NULL)
# ---------------------------------------------------------------------------
$p7 = $d.Paragraphs(22).Range
$xml7 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="009827A8" w:rsidRDefault="009827A8" w:rsidP="009827A8"><w:pPr><w:ind w:left="360"/></w:pPr><w:r w:rsidRPr="009827A8"><w:t>SELECT Maaşı FROM Çalışanlar WHERE Birimi IS NULL</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>
'@
$p7.InsertXML($xml7)

# ---------------------------------------------------------------------------
# Edit 8: "SQL SELECT İfadesi" -> bold
# ---------------------------------------------------------------------------
$p8 = $d.Paragraphs(23).Range
$xml8 = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="008E5471" w:rsidRDefault="008E5471" w:rsidP="008E5471"><w:pPr><w:ind w:left="360"/><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>SQL SELECT İfadesi</w:t></w:r></w:p>
'@
$p8.InsertXML($xml8)
